$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Solar (column E) figures for 2022, 2023, 2024
$ws.Range("E24").Value = 1420
$ws.Range("E25").Value = 1993
$ws.Range("E26").Value = 2127

# Update Energy Storage (column C) figure for 2024
$ws.Range("C26").Value = 32
